$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Rows in the "Ready for handoff" batch generated together
# (2016-09-01 04:24:xx) on the zh-cn / de-de tables and the Overview sheet.
$rows = @(8, 9, 10, 11, 12, 14)

foreach ($r in $rows) {
    # Overview!G<r> = Latest HO Xliff Generate Date (de-de) : 04:24:26 -> 04:24:41
    $wsOverview.Range("G$r").Value = "2016-09-01 04:24:41"

    # zh-cn!H<r> = Latest Handoff Datetime : 04:24:21 -> 04:24:36
    $wsZhCn.Range("H$r").Value = "2016-09-01 04:24:36"

    # de-de!H<r> = Latest Handoff Datetime : 04:24:26 -> 04:24:41
    $wsDeDe.Range("H$r").Value = "2016-09-01 04:24:41"

    # zh-cn!E<r> and de-de!E<r> = Priority : "" -> "ht"
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
